# Updates the cryptos list with freshly scraped values.
# Note: the Price (D) and Volume(1h) (E) columns are stored as TEXT in the
# workbook (e.g. "60.951.61" uses '.' as a thousands separator, which Excel
# would otherwise mis-parse as a number / date). Values that look numeric
# are prefixed with a leading apostrophe so Excel keeps them as text,
# matching the original cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character used in the Pepe price ("0.0<sub>3</sub>0866").
$sub3 = [char]0x2083

# Single-cell price/volume refreshes (no row reshuffle).
$updates = [ordered]@{
    "D2"  = "'60.941.12"
    "E2"  = "  +0.09%  "
    "D3"  = "'2.917.60"
    "E4"  = "  +0.00%  "
    "D5"  = "'590.71"
    "E5"  = "  +0.78%  "
    "D6"  = "'146.40"
    "E6"  = "  +0.70%  "
    "E7"  = "  +0.04%  "
    "E8"  = "  +0.16%  "
    "D9"  = "'6.92"
    "E9"  = "  +1.25%  "
    "E10" = "  -0.73%  "
    "E11" = "  -1.52%  "
    "E12" = "  +0.10%  "
    "E13" = "  +0.00%  "
    "E14" = "  -0.02%  "
    "D15" = "'3.402.70"
    "E15" = "  +0.02%  "
    "D16" = "'60.904.68"
    "E16" = "  +0.06%  "
    "E17" = "  -0.99%  "
    "D18" = "'2.921.58"
    "E18" = "  +0.06%  "
    "D19" = "'432.39"
    "E19" = "  +0.70%  "
    "E20" = "  -1.48%  "
    "E21" = "  -0.53%  "
    "D22" = "'7.11"
    "E22" = "  -0.17%  "
    "D23" = "'81.29"
    "E23" = "  +0.96%  "
    "E24" = "  +1.23%  "
    "E25" = "  -1.35%  "
    "D26" = "'11.85"
    "E26" = "  -1.05%  "
    "E28" = "  +6.34%  "
    "D29" = "'2.61"
    "E29" = "  -0.14%  "
    "E30" = "  -2.80%  "
    "D34" = "'0.0$sub3" + "0866"
    "E34" = "  -0.31%  "
    "E36" = "  -0.23%  "
    "E37" = "  -0.11%  "
    "E38" = "  -1.12%  "
    "E39" = "  -4.77%  "
    "D40" = "'8.56"
    "E40" = "  -1.06%  "
    "D41" = "'41.41"
    "E41" = "  +0.79%  "
    "D42" = "'0.284"
    "E42" = "  -4.20%  "
    "D43" = "'376.47"
    "E43" = "  -0.73%  "
    "D44" = "'2.701.89"
    "E44" = "  -0.01%  "
    "D45" = "'0.0344"
    "E45" = "  -1.84%  "
    "D46" = "'133.68"
    "E46" = "  +0.63%  "
    "D48" = "'23.90"
    "E48" = "  -3.24%  "
    "E49" = "  -0.58%  "
    "E50" = "  -2.73%  "
    "E51" = "  -0.70%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Rows 31/32 swapped rank: Hedera (was rank 29 / row 31) and
# EthereumClassic (was rank 30 / row 32) traded places, each carrying its
# own fresh price/volume data along with it.
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'26.65"
$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.109"
$ws.Range("E32").Value = "  +2.66%  "
